$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column E (Date Sampled) to make room for
# Month / Day / Year columns, shifting Date Sampled..Notes from E:N to H:Q.
$ws.Range("E1:G1").EntireColumn.Insert()

# Match the new columns' width (renders as width 17 without bestFit in the xml).
$ws.Columns("E:G").ColumnWidth = 16.1666666666667

# Headers for the new columns.
$ws.Range("E1").Value = "Month"
$ws.Range("F1").Value = "Day"
$ws.Range("G1").Value = "Year"

# Row 2 - North Head (Date Sampled 7/21/2014)
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 21
$ws.Range("G2").Value = 2014
$ws.Range("M2").Value = 9.1
$ws.Range("N2").Value = 9.6
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 10

# Row 3 - Magic 8 Ball (Date Sampled 7/21/2014)
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 21
$ws.Range("G3").Value = 2014
$ws.Range("M3").Value = 6.3
$ws.Range("N3").Value = 6.7
$ws.Range("O3").Value = 2.5

# Row 4 - Larus Ledge (Date Sampled 7/18/2014)
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 18
$ws.Range("G4").Value = 2014
$ws.Range("M4").Value = 11.4
$ws.Range("N4").Value = 4.4000000000000004
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 12

# Row 5 - Sandpiper Beach (Date Sampled 7/18/2014)
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 18
$ws.Range("G5").Value = 2014
$ws.Range("M5").Value = 9.4
$ws.Range("N5").Value = 4.5
$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 12.5

# Match the author's final selection.
$ws.Range("P4").Select()
